# Timing issue fix - keywords, updated tc1,2 in ubc01
#
# The CasesTab query (cell B2 on the "startup" sheet) is trimmed: the last
# RETURN column (`coalesce(co.cohort_description, '') AS `Cohort``) is
# removed, along with the now-trailing comma on the previous line.
# The SamplesTab (B3) and FilesTab (B4) query cells, and the StatQuery /
# file-name cells (C2:E4), are unchanged in content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$casesTabQuery = @"
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)

MATCH (c)<--(diag:diagnosis)
WHERE s.clinical_study_designation IN ['UBC01'] and demo.breed in ['Mixed Breed', 'Scottish Terrier','Shetland Sheepdog']and diag.disease_term in ['Bladder Cancer','Healthy Control'] and diag.primary_disease_site in ['Bladder', 'Bladder, Prostate', 'Bladder, Urethra', 'Bladder, Urethra, Prostate']
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS ``Case ID`` ,
        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,
        coalesce(s.clinical_study_type, '') AS  ``Study Type``,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,
        coalesce(demo.weight, '') AS ``Weight (kg)``,
        coalesce(diag.best_response, '') AS ``Response to Treatment``
"@

# Here-strings always carry a trailing newline before the closing "@ -- strip it
# so the stored value matches the source cell exactly (no trailing blank line).
$casesTabQuery = $casesTabQuery.TrimEnd("`r", "`n")

$ws.Range("B2").Value2 = $casesTabQuery

# The row shrank by one wrapped line once the "Cohort" column was dropped.
$ws.Rows(2).RowHeight = 304.5

# Matches the author's post-edit selection/viewport in the saved file.
$ws.Range("B2").Select()
